$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (force text format to avoid date-string auto-conversion)
$ws.Range("B1").NumberFormat = "@"
$ws.Range("B1").Value = "http://vocabs.lter-europe.net/so/"
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "https://creativecommons.org/licenses/by/4.0/"
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "1.1.3"
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "2023-08-16"
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "2025-10-16"
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "http://vocabs.lter-europe.net/so/"
$ws.Range("T40").NumberFormat = "@"
$ws.Range("T40").Value = "8/16/2023"
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Eddy Covariance (CO2 flux and concentration, Latent heat flux, Sensible heat flux)"
$ws.Range("U45").NumberFormat = "@"
$ws.Range("U45").Value = "2025-10-13"
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "Vegetation phenology and Leaf Area Index - European scale"
$ws.Range("U47").NumberFormat = "@"
$ws.Range("U47").Value = "2025-10-13"
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Vegetation phenology - site scale"
$ws.Range("U48").NumberFormat = "@"
$ws.Range("U48").Value = "2025-10-13"
$ws.Range("U52").NumberFormat = "@"
$ws.Range("U52").Value = "2025-10-15"
$ws.Range("U53").NumberFormat = "@"
$ws.Range("U53").Value = "2025-10-15"
$ws.Range("B55").NumberFormat = "@"
$ws.Range("B55").Value = "Vegetation aboveground biomass - non-forested sites"
$ws.Range("U55").NumberFormat = "@"
$ws.Range("U55").Value = "2025-10-13"
$ws.Range("B57").NumberFormat = "@"
$ws.Range("B57").Value = "Leaf area index - non-forested sites"
$ws.Range("U57").NumberFormat = "@"
$ws.Range("U57").Value = "2025-10-13"
$ws.Range("B77").NumberFormat = "@"
$ws.Range("B77").Value = "Major ion concentrations: Cl, SO4, Br, Na, K, Mg, Ca, B - groundwater"
$ws.Range("U77").NumberFormat = "@"
$ws.Range("U77").Value = "2025-10-13"
$ws.Range("B83").NumberFormat = "@"
$ws.Range("B83").Value = "Carbon concentration: DOC, DIC - running waters"
$ws.Range("U83").NumberFormat = "@"
$ws.Range("U83").Value = "2025-10-13"
$ws.Range("B85").NumberFormat = "@"
$ws.Range("B85").Value = "Major ion concentrations: Cl, SO4, Br, Na, K, Mg, Ca, B, Silica - running/standing waters"
$ws.Range("U85").NumberFormat = "@"
$ws.Range("U85").Value = "2025-10-13"
$ws.Range("B100").NumberFormat = "@"
$ws.Range("B100").Value = "Land cover, land use, land cover change, land use change (Statistics)"
$ws.Range("U100").NumberFormat = "@"
$ws.Range("U100").Value = "2025-10-13"
$ws.Range("B106").NumberFormat = "@"
$ws.Range("B106").Value = "Resource use"
$ws.Range("U106").NumberFormat = "@"
$ws.Range("U106").Value = "2025-10-13"

# Remove column AG (was entirely empty) to shrink dimension to A1:AF107
$ws.Columns("AG").Delete()

